# Modified problem status in SDE Preparation sheet
#
# 1) Rename the worksheet from "Sheet1" to "Coding Round"
# 2) Mark the 7 "Searching & Sorting" problems in rows 114-120 as solved:
#    - Column B gets the "solved" font style (same style already used on
#      neighbouring rows such as row 101/113 that are marked solved)
#    - Column C switches from "<->" to "Yes"
# 3) Update the sheet's current selection to B124 (matches the author's
#    cursor position after editing row 120 / scrolling further down)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename sheet
$ws.Name = "Coding Round"

# 2) Update rows 114 through 120 (inclusive) in columns B and C.
#    Row 101 is already marked "solved" and has exactly the formatting
#    we want to replicate onto B114:B120 and C114:C120, so copy its
#    per-cell formats across instead of trying to rebuild the font
#    manually.
$doneBSource = $ws.Range("B101")
$doneCSource = $ws.Range("C101")

for ($r = 114; $r -le 120; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $doneBSource.Copy()
    $bCell.PasteSpecial(-4122)   # xlPasteFormats

    $doneCSource.Copy()
    $cCell.PasteSpecial(-4122)   # xlPasteFormats

    $cCell.Value = "Yes"
}

$excel.CutCopyMode = 0

# 3) Update selection / active cell
$null = $ws.Range("B124").Select()
